$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Plano de Desenvolvimento Integrado do Turismo em Minas Gerais" text
$ws.Range("B23").Value = "Plano Diretor do Turismo Verde em Minas Gerais - Plano de Desenvolvimento Integrado do Turismo Sustentável de Minas Gerais"

# 2. Correct the value for row 48 (codigo_iniciativa 9288198)
$ws.Range("D48").Value = 836059044

# 3. Add new row 65: instrumento - deliberação 21
$ws.Range("A58:D58").Copy()
$ws.Range("A65:D65").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A65").Value = 9440688
$ws.Range("B65").Value = "Melhoria da infraestrutura dos municípios – Melhoria da oferta e do acesso à saúde da população de Juiz de Fora e Zona da Mata"
$ws.Range("C65").Value = "IV"
$ws.Range("D65").Value = 150000000
